# Auto-generated edit script: update crypto price/volume table to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Protect numeric-looking Price (column D) values as text, matching original inline-string storage ---
$numericProtectCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D19", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $numericProtectCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# --- Apply updated values ---
$ws.Range("D2").Value = "29.835.67"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").Value = "1.935.70"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "337.23"
$ws.Range("E5").Value = "  +3.72%  "
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "0.4835"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "0.4114"
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("D9").Value = "0.08167"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("D10").Value = "1.014"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").Value = "23.70"
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "6.079"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.902.50"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").Value = "7.278"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "90.88"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "0.06852"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "1.009"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "0.00001035"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").Value = "17.79"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "29.816.42"
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("D22").Value = "5.632"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "11.88"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "2.182"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").Value = "2.183.15"
$ws.Range("E25").Value = "  +3.11%  "
$ws.Range("D26").Value = "6.569"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").Value = "156.93"
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("D28").Value = "20.00"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "2.091"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").Value = "121.17"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D31").Value = "1.008"
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("D32").Value = "0.09652"
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").Value = "5.544"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").Value = "1.412"
$ws.Range("E34").Value = "  +3.53%  "
$ws.Range("D35").Value = "3.533"
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("D36").Value = "0.06589"
$ws.Range("E36").Value = "  +7.78%  "
$ws.Range("D38").Value = "1.201"
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("D39").Value = "0.5985"
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("D40").Value = "10.77"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").Value = "7.954"
$ws.Range("E41").Value = "  -1.14%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "2.486"
$ws.Range("E43").Value = "  +4.27%  "
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").Value = "12.34"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").Value = "0.07481"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("D47").Value = "0.5556"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").Value = "1.985"
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("D49").Value = "116.90"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("D50").Value = "2.416"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").Value = "72.38"
$ws.Range("E51").Value = "  +0.33%  "
